# StageDB.xlsx update:
#   Column J ("StageMapPath" / string / "Prefabs/Stage/CasthleStageMap")
#   is repurposed into a "StageName" (Enum<Sizes>) column whose data rows
#   hold either "Casthle" (rows 4-8) or "Forest" (rows 9-13).
#   Columns G/H (SummonEnemyIDList / SummonBossID) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: column headers (key row)
$ws.Range("J1").Value = "StageName"

# Row 2: column types
$ws.Range("J2").Value = "Enum<Sizes>"

# Row 3: column names (human readable)
$ws.Range("J3").Value = "StageName"

# Rows 4-8 belong to the "Casthle" stage map, rows 9-13 to "Forest"
$ws.Range("J4").Value = "Casthle"
$ws.Range("J5").Value = "Casthle"
$ws.Range("J6").Value = "Casthle"
$ws.Range("J7").Value = "Casthle"
$ws.Range("J8").Value = "Casthle"
$ws.Range("J9").Value = "Forest"
$ws.Range("J10").Value = "Forest"
$ws.Range("J11").Value = "Forest"
$ws.Range("J12").Value = "Forest"
$ws.Range("J13").Value = "Forest"

# The author's selection moved to L4 when they saved the workbook.
$ws.Range("L4").Select() | Out-Null
